$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the organization website from the old domain to the new one.
$ws.Range("B10").Value = "www.stat.gov.kg"

# Leave the selection on the edited cell, matching the saved cursor position.
$ws.Range("B10").Select()
